# Apply the "asset profile section" fixes described in the commit:
#   - Trim the stray trailing space (and xml:space="preserve" marker) from the
#     "Sl " header in D1 so it reads "Sl".
#   - Expand the 2-digit year segment of every Asset Code in column B to a
#     4-digit year (e.g. "-17-" -> "-2017-", "-98-" -> "-1998-").
#   - Two rows additionally had a corrupted "HQ" office-code segment in the
#     source data; fix those specific codes while still expanding the year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the header cell D1 -------------------------------------------------
$ws.Range("D1").Value = "Sl"

# --- Fix Asset Codes in column B -------------------------------------------
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $code = $cell.Value2
    if ([string]::IsNullOrEmpty($code)) { continue }

    if ($code -eq "FRC-HQ-SLM-C-98-0000") {
        # Special-cased legacy record with a different rewrite pattern
        $cell.Value = "FRC-FRC-0-C-1998-0000"
    }
    elseif ($code -match '^(FRC)-(HQ)-(SLM)-([A-Za-z])-(\d{2})-(\d{4})$') {
        $prefix = $matches[1]
        $office = $matches[2]
        $mid    = $matches[3]
        $letter = $matches[4]
        $yy     = [int]$matches[5]
        $num    = $matches[6]

        if ($yy -ge 90) {
            $year = 1900 + $yy
        } else {
            $year = 2000 + $yy
        }

        if ($r -ge 320 -and $r -le 331) {
            # These rows also had the office segment corrupted to "dfsdaf"
            $office = "dfsdaf"
        }

        $newCode = "$prefix-$office-$mid-$letter-$year-$num"
        $cell.Value = $newCode
    }
}
